$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.006.63'
$ws.Range('D3').Value = '1.831.08'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('D4').Value = "'0.9992"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'241.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = "'0.6269"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.42%  '
$ws.Range('D7').Value = "'1.001"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = "'0.07601"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.27%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = "'0.2913"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = "'22.72"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.10%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'0.07748"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.831.81'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.955"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.6627"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'82.42"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = "'0.000009461"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +9.54%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = "'5.983"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.82%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '29.002.65'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'223.95"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = "'12.33"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = "'1.0000"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = "'7.210"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = "'1.001"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = "'158.76"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = "'8.410"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = "'0.1365"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.40%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'17.83"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = "'1.489"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.40%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = "'4.055"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'4.023"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = "'1.192"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.05190"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = "'1.844"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'0.7376"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'1.145"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.698"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.253.17'
$ws.Range('E37').Value = '  -4.77%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = "'2.761"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.01788"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'6.211"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.31%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'0.8929"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.63%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'1.001"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = "'101.46"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.975.85'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'64.53"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = "'0.5109"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = "'0.00000000120"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').Value = "'0.3977"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'8.838"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.05751"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = "'6.680"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.89%  '
